$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "23.700.20"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "1.654.07"
$ws.Range("E3").Value = "  +1.12%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("E5").Value = "  +0.04%  "
Set-TextValue $ws.Range("D6") "303.10"
$ws.Range("E6").Value = "  -0.17%  "
Set-TextValue $ws.Range("D7") "0.3802"
$ws.Range("E7").Value = "  +0.41%  "
Set-TextValue $ws.Range("D8") "0.3620"
Set-TextValue $ws.Range("D9") "51.23"
$ws.Range("E9").Value = "  -1.00%  "
Set-TextValue $ws.Range("D10") "1.245"
$ws.Range("E10").Value = "  +1.31%  "
Set-TextValue $ws.Range("D11") "0.08215"
$ws.Range("E11").Value = "  +0.39%  "
Set-TextValue $ws.Range("D12") "1.002"
$ws.Range("E12").Value = "  +0.18%  "
Set-TextValue $ws.Range("D13") "22.63"
$ws.Range("E13").Value = "  +1.09%  "
Set-TextValue $ws.Range("D14") "6.522"
$ws.Range("E14").Value = "  +0.76%  "
Set-TextValue $ws.Range("D15") "7.439"
$ws.Range("E15").Value = "  +0.88%  "
Set-TextValue $ws.Range("D16") "0.00001234"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "1.659.38"
$ws.Range("E17").Value = "  +1.90%  "
Set-TextValue $ws.Range("D18") "97.36"
$ws.Range("E18").Value = "  +2.33%  "
Set-TextValue $ws.Range("D19") "0.07021"
$ws.Range("E19").Value = "  +1.09%  "
Set-TextValue $ws.Range("D20") "6.805"
$ws.Range("E20").Value = "  +3.33%  "
Set-TextValue $ws.Range("D21") "17.70"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("E22").Value = "  +0.02%  "
Set-TextValue $ws.Range("D23") "12.87"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D24").Value = "23.711.90"
$ws.Range("E24").Value = "  +1.14%  "
Set-TextValue $ws.Range("D25") "2.519"
$ws.Range("E25").Value = "  +0.33%  "
Set-TextValue $ws.Range("D26") "3.046"
$ws.Range("E26").Value = "  -0.30%  "
Set-TextValue $ws.Range("D27") "21.28"
$ws.Range("E27").Value = "  +0.65%  "
Set-TextValue $ws.Range("D28") "153.43"
$ws.Range("E28").Value = "  +1.68%  "
Set-TextValue $ws.Range("D29") "5.236"
$ws.Range("E29").Value = "  -0.71%  "
Set-TextValue $ws.Range("D30") "134.62"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "1.839.58"
$ws.Range("E31").Value = "  +1.58%  "
Set-TextValue $ws.Range("D32") "6.931"
$ws.Range("E32").Value = "  +4.85%  "
Set-TextValue $ws.Range("D33") "2.223"
$ws.Range("E33").Value = "  +2.64%  "
Set-TextValue $ws.Range("D34") "1.069"
$ws.Range("E34").Value = "  +2.24%  "
Set-TextValue $ws.Range("D35") "11.70"
$ws.Range("E35").Value = "  +4.20%  "
Set-TextValue $ws.Range("D36") "0.02814"
$ws.Range("E36").Value = "  +2.27%  "
Set-TextValue $ws.Range("D37") "0.2533"
$ws.Range("E37").Value = "  +1.58%  "
Set-TextValue $ws.Range("D38") "0.08800"
$ws.Range("E38").Value = "  +0.26%  "
Set-TextValue $ws.Range("D39") "6.100"
$ws.Range("E39").Value = "  +1.23%  "
Set-TextValue $ws.Range("D41") "12.97"
$ws.Range("E41").Value = "  +6.78%  "
Set-TextValue $ws.Range("D42") "0.7029"
$ws.Range("E42").Value = "  +0.38%  "
Set-TextValue $ws.Range("D43") "1.336"
$ws.Range("E43").Value = "  -0.32%  "
Set-TextValue $ws.Range("D44") "16.00"
$ws.Range("E44").Value = "  +1.00%  "
Set-TextValue $ws.Range("D45") "0.6506"
$ws.Range("E45").Value = "  +0.09%  "
Set-TextValue $ws.Range("D46") "2.315"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("E47").Value = "  +0.01%  "
Set-TextValue $ws.Range("D48") "3.979"
$ws.Range("E48").Value = "  +0.25%  "
Set-TextValue $ws.Range("D49") "0.07961"
$ws.Range("E49").Value = "  -0.23%  "
Set-TextValue $ws.Range("D50") "128.05"
$ws.Range("E50").Value = "  +0.66%  "
Set-TextValue $ws.Range("D51") "1.191"
$ws.Range("E51").Value = "  -0.02%  "
